$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder comma-separated control lists in column A (IA Control) ---
# Values are unchanged sets of tokens; only the order changes, matching the
# upstream srg_mapping content regeneration.
$ws.Range('A2').Value2 = 'AU-4 (1),AU-4'
$ws.Range('A3').Value2 = 'AU-14 (1),AU-4'
$ws.Range('A6').Value2 = 'AU-12 (3),AU-7 a,AC-6 (9),AU-8 b,CM-5 (1),AU-7 b,AC-6 (8)'
$ws.Range('A7').Value2 = 'AU-12 (3),AU-12 c,CM-6 b,AU-7 a,AU-8 b,AU-12 a,CM-5 (1),AU-7 b'
$ws.Range('A11').Value2 = 'IA-2 (11),IA-2 (12)'
$ws.Range('A12').Value2 = 'IA-2 (11),IA-2 (12)'
$ws.Range('A14').Value2 = 'CM-7 (2),CM-7 (5) (b)'
$ws.Range('A15').Value2 = 'CM-7 (2),CM-7 (5) (b)'
$ws.Range('A17').Value2 = 'CM-7 (2),CM-6 b'
$ws.Range('A22').Value2 = 'CM-7 (2),CM-6 b'
$ws.Range('A23').Value2 = 'CM-7 (2),CM-6 b'
$ws.Range('A38').Value2 = 'AC-7 a,AC-7 b'
$ws.Range('A39').Value2 = 'AC-7 a,AC-7 b'
$ws.Range('A40').Value2 = 'AC-7 a,AC-7 b'
$ws.Range('A41').Value2 = 'AC-7 a,AC-7 b'
$ws.Range('A45').Value2 = 'AU-3 (1),IA-2,IA-8'
$ws.Range('A46').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A47').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A48').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A49').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A50').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A51').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A52').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A53').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A54').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A55').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A56').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A57').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A58').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A59').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A60').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A61').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A62').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A63').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A64').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A65').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A66').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A67').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A68').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A69').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A70').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A71').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A72').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A73').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A74').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A75').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A76').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A77').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A78').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A79').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A80').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A81').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A82').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A83').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A84').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A85').Value2 = 'AU-3,MA-4 (1) (a),AU-3 (1)'
$ws.Range('A86').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A87').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A88').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A89').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A90').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A91').Value2 = 'MA-4 (1) (a),AU-3 (1),AU-12 c'
$ws.Range('A92').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A93').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A94').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A95').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A96').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A97').Value2 = 'AU-3,MA-4 (1) (a),AU-3 (1),AU-12 c'
$ws.Range('A98').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A99').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A100').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A101').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A102').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A103').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A104').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A105').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A106').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-3,MA-4 (1) (a)'
$ws.Range('A107').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A108').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A109').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A110').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A111').Value2 = 'AU-12 c,AC-2 (4),AU-3 (1),AU-12 a,AU-3,MA-4 (1) (a)'
$ws.Range('A112').Value2 = 'AU-12 c,AU-3 (1),AU-12 a,AU-3,AU-14 (1),MA-4 (1) (a)'
$ws.Range('A113').Value2 = 'AC-6 (10),CM-6 b'
$ws.Range('A114').Value2 = 'AC-6 (10),CM-6 b'
$ws.Range('A115').Value2 = 'AC-6 (10),AC-11 b'
$ws.Range('A116').Value2 = 'AC-6 (10),CM-6 b'
$ws.Range('A120').Value2 = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c'
$ws.Range('A121').Value2 = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c'
$ws.Range('A122').Value2 = 'AU-9,AU-12 c'
$ws.Range('A127').Value2 = 'CM-5 (1),AU-12 c,AC-2 (4),AC-6 (9)'
$ws.Range('A129').Value2 = 'IA-5 (1) (a),IA-5 (1) (b),CM-6 b'
$ws.Range('A133').Value2 = 'SC-13,AC-17 (2),MA-4 c,SC-8'
$ws.Range('A134').Value2 = 'SC-10,AC-12,MA-4 e,MA-4 (7)'
$ws.Range('A135').Value2 = 'SC-10,AC-12'
$ws.Range('A136').Value2 = 'SC-10,AC-12'
$ws.Range('A138').Value2 = 'AU-7 (1),MA-4 (1) (a),CM-6 b,AU-7 a,AU-3 (1),AU-12 a,AU-3,CM-5 (1),AU-14 (1),AU-6 (4)'
$ws.Range('A143').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A144').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A145').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A146').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A147').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A148').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A149').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A150').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A151').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A152').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A153').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A154').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A155').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A166').Value2 = 'SC-8 (2),SC-8 (1),SC-8'
$ws.Range('A168').Value2 = 'AC-17 (2),SC-8'
$ws.Range('A169').Value2 = 'SC-8 (2),SC-8'
$ws.Range('A172').Value2 = 'AC-11 a,AC-11 b'
$ws.Range('A173').Value2 = 'AC-11 a,AC-11 b'
$ws.Range('A174').Value2 = 'AC-11 a,AC-11 b'
$ws.Range('A179').Value2 = 'AU-4 (1),AU-6 (4),CM-6 b'
$ws.Range('A180').Value2 = 'CM-7 b,AC-17 (9),AC-17 (1),CM-6 b'
$ws.Range('A181').Value2 = 'CM-7 b,AC-17 (1),CM-6 b'
$ws.Range('A193').Value2 = 'AU-3,CM-6 b'
$ws.Range('A199').Value2 = 'AU-4 (1),AU-3'
$ws.Range('A206').Value2 = 'AU-4 (1),CM-6 b'
$ws.Range('A207').Value2 = 'SC-28 (1),SC-28'
$ws.Range('A211').Value2 = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range('A212').Value2 = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range('A213').Value2 = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range('A214').Value2 = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range('A215').Value2 = 'AU-12 c,AC-2 (4),AC-6 (9)'
$ws.Range('A220').Value2 = 'IA-2 (5),CM-6 b'
$ws.Range('A221').Value2 = 'IA-2 (5),IA-2 (2),IA-2,IA-2 (3),IA-2 (4)'
$ws.Range('A222').Value2 = 'IA-2 (5),IA-2 (2),IA-2,IA-2 (3),IA-2 (4)'
$ws.Range('A223').Value2 = 'SC-8 (1),AC-18 (1),SC-8'
$ws.Range('A225').Value2 = 'IA-7,IA-5 (1) (c)'
$ws.Range('A226').Value2 = 'IA-7,CM-6 b'
$ws.Range('A227').Value2 = 'IA-7,CM-6 b'
$ws.Range('A228').Value2 = 'IA-7,CM-6 b'
$ws.Range('A230').Value2 = 'CM-7 a,IA-7'
$ws.Range('A241').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A242').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A244').Value2 = 'SC-2,SI-16,CM-6 b'
$ws.Range('A246').Value2 = 'SC-3,SI-16'
$ws.Range('A269').Value2 = 'CM-6 b,IA-2 (2)'
$ws.Range('A270').Value2 = 'IA-2 (3),IA-2 (4),IA-2 (1),IA-2 (2)'
$ws.Range('A272').Value2 = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A275').Value2 = 'SC-4,CM-6 b'
$ws.Range('A276').Value2 = 'SC-4,SC-2'
$ws.Range('A277').Value2 = 'SC-4,SC-2'
$ws.Range('A297').Value2 = 'IA-2 (11),IA-2 (12)'
$ws.Range('A299').Value2 = 'IA-2 (1),IA-2 (11),IA-2 (12)'
$ws.Range('A309').Value2 = 'AU-8 b,AU-8 (1) (a),AU-8 (1) (b)'
$ws.Range('A328').Value2 = 'CM-5 (1),AU-12 c'
$ws.Range('A342').Value2 = 'CM-7 b,IA-3'
$ws.Range('A343').Value2 = 'CM-7 a,CM-7 b'
$ws.Range('A344').Value2 = 'CM-7 a,CM-7 b'
$ws.Range('A345').Value2 = 'CM-7 b,AC-17 (1)'
$ws.Range('A346').Value2 = 'CM-7 a,AC-18 (1)'
$ws.Range('A347').Value2 = 'CM-7 a,IA-5 (1) (c),CM-6 b'
$ws.Range('A361').Value2 = 'SI-6 b,CM-3 (5),SI-6 d'
$ws.Range('A362').Value2 = 'CM-7 a,CM-6 b'
$ws.Range('A367').Value2 = 'CM-7 a,SI-16'
$ws.Range('A374').Value2 = 'CM-7 a,CM-6 b'
$ws.Range('A375').Value2 = 'CM-7 a,CM-6 b'
$ws.Range('A376').Value2 = 'CM-7 a,CM-6 b'
$ws.Range('A391').Value2 = 'IA-5 (1) (a),CM-6 b'
$ws.Range('A397').Value2 = 'CM-3 (5),SI-6 d'
$ws.Range('A401').Value2 = 'SC-3,CM-6 b'
$ws.Range('A402').Value2 = 'SC-3,CM-6 b'
$ws.Range('A403').Value2 = 'SC-3,CM-6 b'
$ws.Range('A450').Value2 = 'CM-5 (1),CM-6 b'
$ws.Range('A451').Value2 = 'CM-5 (1),CM-6 b'

# --- Update H466 / K466 text to mention recovery boot / GRUB_DISABLE_RECOVERY ---
$hCell = $ws.Range('H466')
$hFind = 'Using interactive boot, the console user could disable auditing, firewalls,'
$hReplace = 'Using interactive or recovery boot, the console user could disable auditing, firewalls,'
$hCell.Value2 = $hCell.Value2.Replace($hFind, $hReplace)

$kCell = $ws.Range('K466')
$kFind = "that interactive boot is enabled at boot time.`n`nIf Interactive boot is enabled at boot time then this is a finding."
$kReplace = "that interactive boot is enabled at boot time and verify that`n GRUB_DISABLE_RECOVERY=true  to disable recovery boot.`nIf Interactive boot is enabled at boot time then this is a finding."
$kCell.Value2 = $kCell.Value2.Replace($kFind, $kReplace)

Write-Host "Done applying srg-mapping-rhel9 control reorderings and recovery-boot text update."
